# AFDP-9108: Fix change-consultation-status workflow rule so the approvers /
# reviewers modal shows up correctly.
#
# The "Form Workflow Rules" table (Sheet1) drives which Drools rule fires
# for a given ACM file type. Two problems were found with the rule table:
#
#   1. Column D ("Start a Workflow Process? (true/false)") was populated
#      with the literal text "true"/"false" instead of real boolean values,
#      which the rule engine does not evaluate as a boolean.
#   2. The "Change Consultation Status" row (row 24) re-used the
#      "change_case_status" file-type key from the Case rule above it
#      instead of its own "change_consultation_status" key, so the
#      consultation-specific rule never matched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 = "Default Workflow (no other rules match)" -> should not start a
# workflow process.
$ws.Range("D17").Value = $false

# Rows 18-24 = the specific rules -> each of them does start a workflow
# process, now stored as a true boolean instead of the text "true"/"true ".
$ws.Range("D18:D24").Value = $true

# Row 24 = "Change Consultation Status" rule must use its own file type key.
$ws.Range("C24").Value = "change_consultation_status"
